$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "y" for two more rows (Search for user nearby row, and signup bug row)
$ws.Range("C6").Value = "y"
$ws.Range("C17").Value = "y"

# Update the last active selection to C16
$ws.Range("C16").Select()
